# Bump the published "term" package to 1.1.0:
#  - Metadata!B3 (Version)   1.0.0 -> 1.1.0
#  - Metadata!B8 (Date)      2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item(1)

$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
